$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12
$ws.Range("A12").Value = "Ports USB "
$ws.Range("B12").Value = [datetime]"2019-02-07"
$ws.Range("C12").Value = "3h00"
$ws.Range("D12").Value = "Je n'ai définitvement trouvé aucun résultat pour récupérer les emplacements des clés en c#"

# Row 13
$ws.Range("A13").Value = "Ports USB"
$ws.Range("B13").Value = [datetime]"2019-02-08"
$ws.Range("C13").Value = "2h00"
$ws.Range("D13").Value = "J'ai décidé d'utiliser powerShell pour récupérer les ports des clés, j'ai trouvé un moyen d'exécuter un script en C# et d'en récupérer les informations"

# Row 14
$ws.Range("A14").Value = "Debugage"
$ws.Range("B14").Value = [datetime]"2019-02-08"
$ws.Range("C14").Value = "1h00"
$ws.Range("D14").Value = "Débugage de quelques bugs notamment en ce qui concernait la mise à jour des clés existantes"

$ws.Range("D14").Select() | Out-Null
